# Auto-generated PowerShell Excel COM-interop script
# Applies numeric cell value updates to the Seraph_Profits leve-profit tables
# across all 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (40 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 323.5
$ws.Range("J33").Value = 1111.25
$ws.Range("L33").Value = 1111.25
$ws.Range("N33").Value = -1569.25
$ws.Range("H107").Value = 391.8421
$ws.Range("I107").Value = 402.5
$ws.Range("K107").Value = 402.5
$ws.Range("M107").Value = 1517.5
$ws.Range("H132").Value = 2881.9
$ws.Range("I132").Value = 2915.375
$ws.Range("J132").Value = 2748
$ws.Range("K132").Value = 8746.125
$ws.Range("L132").Value = 8244
$ws.Range("M132").Value = -6216.125
$ws.Range("N132").Value = -13304
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 2623.7666
$ws.Range("I137").Value = 1252.4166
$ws.Range("J137").Value = 3538
$ws.Range("K137").Value = 3757.2498
$ws.Range("L137").Value = 10614
$ws.Range("M137").Value = -1207.2498
$ws.Range("N137").Value = -15714
$ws.Range("H138").Value = 7649.478
$ws.Range("I138").Value = 4578.6
$ws.Range("J138").Value = 8502.5
$ws.Range("K138").Value = 13735.8
$ws.Range("L138").Value = 25507.5
$ws.Range("M138").Value = -8595.800000000001
$ws.Range("N138").Value = -35787.5
$ws.Range("H141").Value = 11364
$ws.Range("I141").Value = 11364
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 34092
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -28912
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM (51 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1093.091
$ws.Range("I2").Value = 1002.7778
$ws.Range("K2").Value = 1002.7778
$ws.Range("M2").Value = -889.7778
$ws.Range("H45").Value = 1815.5
$ws.Range("I45").Value = 1795
$ws.Range("K45").Value = 1795
$ws.Range("M45").Value = -1418
$ws.Range("H61").Value = 1668.2354
$ws.Range("I61").Value = 1600.6666
$ws.Range("J61").Value = 2175
$ws.Range("K61").Value = 1600.6666
$ws.Range("L61").Value = 2175
$ws.Range("M61").Value = -1388.6666
$ws.Range("N61").Value = -2599
$ws.Range("H74").Value = 2135.1177
$ws.Range("I74").Value = 1001.6667
$ws.Range("J74").Value = 3410.25
$ws.Range("K74").Value = 1001.6667
$ws.Range("L74").Value = 3410.25
$ws.Range("M74").Value = -127.6667
$ws.Range("N74").Value = -5158.25
$ws.Range("H77").Value = 2135.1177
$ws.Range("I77").Value = 1001.6667
$ws.Range("J77").Value = 3410.25
$ws.Range("K77").Value = 5008.3335
$ws.Range("L77").Value = 17051.25
$ws.Range("M77").Value = -640.3334999999997
$ws.Range("N77").Value = -25787.25
$ws.Range("H116").Value = 1093.091
$ws.Range("I116").Value = 1002.7778
$ws.Range("K116").Value = 1002.7778
$ws.Range("M116").Value = 1291.2222
$ws.Range("H122").Value = 3765.9
$ws.Range("I122").Value = 3272.7144
$ws.Range("K122").Value = 9818.143199999999
$ws.Range("M122").Value = -7368.143199999999
$ws.Range("H132").Value = 1894.2094
$ws.Range("I132").Value = 1683.5641
$ws.Range("J132").Value = 3948
$ws.Range("K132").Value = 5050.692300000001
$ws.Range("L132").Value = 11844
$ws.Range("M132").Value = -2520.692300000001
$ws.Range("N132").Value = -16904
$ws.Range("H136").Value = 1668.2354
$ws.Range("I136").Value = 1600.6666
$ws.Range("J136").Value = 2175
$ws.Range("K136").Value = 4801.9998
$ws.Range("L136").Value = 6525
$ws.Range("M136").Value = -2251.9998
$ws.Range("N136").Value = -11625

# ---- Sheet: BSM (12 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1093.091
$ws.Range("I3").Value = 1002.7778
$ws.Range("K3").Value = 1002.7778
$ws.Range("M3").Value = -888.7778
$ws.Range("H64").Value = 1044.6923
$ws.Range("J64").Value = 1193.5714
$ws.Range("L64").Value = 1193.5714
$ws.Range("N64").Value = -1643.5714
$ws.Range("H67").Value = 1044.6923
$ws.Range("J67").Value = 1193.5714
$ws.Range("L67").Value = 1193.5714
$ws.Range("N67").Value = -2753.5714

# ---- Sheet: CRP (4 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 52500
$ws.Range("J96").Value = 52500
$ws.Range("L96").Value = 52500
$ws.Range("N96").Value = -57992

# ---- Sheet: CUL (47 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1500
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4388
$ws.Range("N3").ClearContents()
$ws.Range("H23").Value = 256.25
$ws.Range("J23").Value = 250
$ws.Range("L23").Value = 750
$ws.Range("N23").Value = -1220
$ws.Range("H113").Value = 1883.0769
$ws.Range("I113").Value = 2333
$ws.Range("K113").Value = 6999
$ws.Range("M113").Value = -4829
$ws.Range("H130").Value = 5580
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 6475
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 19425
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -29465
$ws.Range("H132").Value = 3332
$ws.Range("I132").Value = 998.1667
$ws.Range("J132").Value = 7999.6665
$ws.Range("K132").Value = 8983.5003
$ws.Range("L132").Value = 71996.9985
$ws.Range("M132").Value = -6453.5003
$ws.Range("N132").Value = -77056.9985
$ws.Range("H134").Value = 3010.125
$ws.Range("I134").Value = 1297.4286
$ws.Range("K134").Value = 3892.2858
$ws.Range("M134").Value = 1177.7142
$ws.Range("H136").Value = 12795
$ws.Range("I136").Value = 12795
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 38385
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -33285
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 3285.652
$ws.Range("I139").Value = 2269
$ws.Range("J139").Value = 6166.1665
$ws.Range("K139").Value = 6807
$ws.Range("L139").Value = 18498.4995
$ws.Range("M139").Value = -1667
$ws.Range("N139").Value = -28778.4995

# ---- Sheet: GSM (4 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 613836.9399999999
$ws.Range("I122").Value = 78852.234
$ws.Range("K122").Value = 236556.702
$ws.Range("M122").Value = -234106.702

# ---- Sheet: LTW (69 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8162.933
$ws.Range("I16").Value = 6287.1665
$ws.Range("K16").Value = 6287.1665
$ws.Range("M16").Value = -6117.1665
$ws.Range("H22").Value = 814.8
$ws.Range("I22").Value = 843.5
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 843.5
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -548.5
$ws.Range("N22").Value = -1290
$ws.Range("H27").Value = 814.8
$ws.Range("I27").Value = 843.5
$ws.Range("J27").Value = 700
$ws.Range("K27").Value = 843.5
$ws.Range("L27").Value = 700
$ws.Range("M27").Value = -736.5
$ws.Range("N27").Value = -914
$ws.Range("H40").Value = 2272.7273
$ws.Range("I40").Value = 1858.2858
$ws.Range("J40").Value = 2998
$ws.Range("K40").Value = 1858.2858
$ws.Range("L40").Value = 2998
$ws.Range("M40").Value = -1722.2858
$ws.Range("N40").Value = -3270
$ws.Range("H46").Value = 4090.9092
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312
$ws.Range("H55").Value = 386.14285
$ws.Range("I55").Value = 364.86667
$ws.Range("J55").Value = 439.33334
$ws.Range("K55").Value = 364.86667
$ws.Range("L55").Value = 439.33334
$ws.Range("M55").Value = -191.86667
$ws.Range("N55").Value = -785.33334
$ws.Range("H61").Value = 6248.143
$ws.Range("I61").Value = 7147.2
$ws.Range("J61").Value = 4000.5
$ws.Range("K61").Value = 7147.2
$ws.Range("L61").Value = 4000.5
$ws.Range("M61").Value = -6945.2
$ws.Range("N61").Value = -4404.5
$ws.Range("H74").Value = 42999
$ws.Range("I74").Value = 42999
$ws.Range("K74").Value = 42999
$ws.Range("M74").Value = -42001
$ws.Range("H77").Value = 42999
$ws.Range("I77").Value = 42999
$ws.Range("K77").Value = 128997
$ws.Range("M77").Value = -124005
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H113").Value = 6248.143
$ws.Range("I113").Value = 7147.2
$ws.Range("J113").Value = 4000.5
$ws.Range("K113").Value = 7147.2
$ws.Range("L113").Value = 4000.5
$ws.Range("M113").Value = -4977.2
$ws.Range("N113").Value = -8340.5
$ws.Range("H122").Value = 3524.5
$ws.Range("I122").Value = 3500.9
$ws.Range("J122").Value = 3642.5
$ws.Range("K122").Value = 10502.7
$ws.Range("L122").Value = 10927.5
$ws.Range("M122").Value = -8052.700000000001
$ws.Range("N122").Value = -15827.5

# ---- Sheet: WVR (11 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1320.7858
$ws.Range("I107").Value = 661.875
$ws.Range("J107").Value = 2199.3333
$ws.Range("K107").Value = 1985.625
$ws.Range("L107").Value = 6597.999899999999
$ws.Range("M107").Value = -65.625
$ws.Range("N107").Value = -10437.9999
$ws.Range("H132").Value = 2204.7407
$ws.Range("J132").Value = 3076.375
$ws.Range("L132").Value = 9229.125
$ws.Range("N132").Value = -14289.125

